$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3991.6667
$ws.Range("J51").Value = 4389.8
$ws.Range("L51").Value = 4389.8
$ws.Range("N51").Value = -5357.8

$ws.Range("H76").Value = 83423110
$ws.Range("I76").Value = 205880.4
$ws.Range("K76").Value = 205880.4
$ws.Range("M76").Value = -205565.4

$ws.Range("H79").Value = 83423110
$ws.Range("I79").Value = 205880.4
$ws.Range("K79").Value = 205880.4
$ws.Range("M79").Value = -204788.4

$ws.Range("H112").Value = 1955.625
$ws.Range("J112").Value = 1963.9286
$ws.Range("L112").Value = 5891.7858
$ws.Range("N112").Value = -8107.7858

$ws.Range("H118").Value = 32231.5
$ws.Range("I118").Value = 38587.8
$ws.Range("K118").Value = 115763.4
$ws.Range("M118").Value = -114106.4

$ws.Range("H127").Value = 3237.0435
$ws.Range("I127").Value = 950.2727
$ws.Range("J127").Value = 5333.25
$ws.Range("K127").Value = 2850.8181
$ws.Range("L127").Value = 15999.75
$ws.Range("M127").Value = 2109.1819
$ws.Range("N127").Value = -25919.75

$ws.Range("H138").Value = 5484.2144
$ws.Range("J138").Value = 6108.273
$ws.Range("L138").Value = 18324.819
$ws.Range("N138").Value = -28604.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 830.4194
$ws.Range("I32").Value = 823.9551
$ws.Range("K32").Value = 823.9551
$ws.Range("M32").Value = -536.9551

$ws.Range("H45").Value = 64642.062
$ws.Range("I45").Value = 73305.21000000001
$ws.Range("K45").Value = 73305.21000000001
$ws.Range("M45").Value = -72928.21000000001

$ws.Range("H125").Value = 56665.832
$ws.Range("J125").Value = 56665.832
$ws.Range("L125").Value = 56665.832
$ws.Range("N125").Value = -66505.83199999999

$ws.Range("H133").Value = 58571.43
$ws.Range("J133").Value = 58571.43
$ws.Range("L133").Value = 58571.43
$ws.Range("N133").Value = -63631.43

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2517
$ws.Range("I20").Value = 2004.1666
$ws.Range("J20").Value = 3029.8333
$ws.Range("K20").Value = 2004.1666
$ws.Range("L20").Value = 3029.8333
$ws.Range("M20").Value = -1757.1666
$ws.Range("N20").Value = -3523.8333

$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("N45").ClearContents()

$ws.Range("H134").Value = 43115.535
$ws.Range("I134").Value = 2848.5715
$ws.Range("K134").Value = 8545.7145
$ws.Range("M134").Value = -6010.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2259.182
$ws.Range("I122").Value = 1793.2
$ws.Range("K122").Value = 5379.6
$ws.Range("M122").Value = -2929.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16254838
$ws.Range("I4").Value = 945535.5
$ws.Range("K4").Value = 2836606.5
$ws.Range("M4").Value = -2836494.5

$ws.Range("H114").Value = 697.2222
$ws.Range("J114").Value = 1007.5
$ws.Range("L114").Value = 3022.5
$ws.Range("N114").Value = -9530.5

$ws.Range("H117").Value = 516.6667
$ws.Range("I117").Value = 250
$ws.Range("J117").Value = 650
$ws.Range("K117").Value = 750
$ws.Range("L117").Value = 1950
$ws.Range("M117").Value = 2692
$ws.Range("N117").Value = -8834

$ws.Range("H128").Value = 148285.14
$ws.Range("I128").Value = 148285.14
$ws.Range("K128").Value = 444855.42
$ws.Range("M128").Value = -439875.42

$ws.Range("H131").Value = 58442.082
$ws.Range("J131").Value = 37752.69
$ws.Range("L131").Value = 113258.07
$ws.Range("N131").Value = -123338.07

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 3975
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 3975
$ws.Range("K4").Value = 0
$ws.Range("N4").Value = -4199
$ws.Range("M4").ClearContents()

$ws.Range("H70").Value = 7620
$ws.Range("I70").Value = 7292.3335
$ws.Range("J70").Value = 8209.799999999999
$ws.Range("K70").Value = 7292.3335
$ws.Range("L70").Value = 8209.799999999999
$ws.Range("M70").Value = -7022.3335
$ws.Range("N70").Value = -8749.799999999999

$ws.Range("H73").Value = 7620
$ws.Range("I73").Value = 7292.3335
$ws.Range("J73").Value = 8209.799999999999
$ws.Range("K73").Value = 7292.3335
$ws.Range("L73").Value = 8209.799999999999
$ws.Range("M73").Value = -6356.3335
$ws.Range("N73").Value = -10081.8

$ws.Range("H120").Value = 37033.5
$ws.Range("J120").Value = 37033.5
$ws.Range("L120").Value = 37033.5
$ws.Range("N120").Value = -46709.5

$ws.Range("H122").Value = 653444.5600000001
$ws.Range("I122").Value = 791827.5
$ws.Range("K122").Value = 2375482.5
$ws.Range("M122").Value = -2373032.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 594434.0600000001
$ws.Range("I40").Value = 672931.9399999999
$ws.Range("K40").Value = 672931.9399999999
$ws.Range("M40").Value = -672795.9399999999

$ws.Range("H61").Value = 3823.9656
$ws.Range("I61").Value = 2849.05
$ws.Range("K61").Value = 2849.05
$ws.Range("M61").Value = -2647.05

$ws.Range("H68").Value = 75009
$ws.Range("I68").Value = 3374.5
$ws.Range("J68").Value = 170521.67
$ws.Range("K68").Value = 3374.5
$ws.Range("L68").Value = 170521.67
$ws.Range("M68").Value = -2625.5
$ws.Range("N68").Value = -172019.67

$ws.Range("H71").Value = 75009
$ws.Range("I71").Value = 3374.5
$ws.Range("J71").Value = 170521.67
$ws.Range("K71").Value = 16872.5
$ws.Range("L71").Value = 852608.3500000001
$ws.Range("M71").Value = -13128.5
$ws.Range("N71").Value = -860096.3500000001

$ws.Range("H113").Value = 3823.9656
$ws.Range("I113").Value = 2849.05
$ws.Range("K113").Value = 2849.05
$ws.Range("M113").Value = -679.0500000000002

$ws.Range("H121").Value = 34169
$ws.Range("J121").Value = 34169
$ws.Range("L121").Value = 34169
$ws.Range("N121").Value = -37663

$ws.Range("H122").Value = 616281.4
$ws.Range("I122").Value = 4010
$ws.Range("J122").Value = 1106098.5
$ws.Range("K122").Value = 12030
$ws.Range("L122").Value = 3318295.5
$ws.Range("M122").Value = -9580
$ws.Range("N122").Value = -3323195.5

$ws.Range("H136").Value = 8015400
$ws.Range("I136").Value = 16020086
$ws.Range("K136").Value = 48060258
$ws.Range("M136").Value = -48057708

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7293.2666
$ws.Range("J62").Value = 7999.8335
$ws.Range("L62").Value = 7999.8335
$ws.Range("N62").Value = -9247.833500000001

$ws.Range("H65").Value = 7293.2666
$ws.Range("J65").Value = 7999.8335
$ws.Range("L65").Value = 39999.1675
$ws.Range("N65").Value = -46239.1675

$ws.Range("H95").Value = 79998
$ws.Range("J95").Value = 79998
$ws.Range("L95").Value = 79998
$ws.Range("N95").Value = -85490

$ws.Range("H113").Value = 2038.1538
$ws.Range("I113").Value = 2110
$ws.Range("J113").Value = 1876.5
$ws.Range("K113").Value = 6330
$ws.Range("L113").Value = 5629.5
$ws.Range("M113").Value = -4160
$ws.Range("N113").Value = -9969.5

$ws.Range("H122").Value = 3680.2917
$ws.Range("I122").Value = 2964.6316
$ws.Range("K122").Value = 8893.8948
$ws.Range("M122").Value = -6443.8948
